$d = $word.ActiveDocument

# --- Step 1: insert a new, empty paragraph right after the paragraph that
# ends with "Once you've changed follow the next step in your assignment
# task." -- it automatically inherits the same paragraph border (pBdr)
# formatting via the Find/Replace paragraph-mark insertion.
$anchorText = "Once you" + [char]8217 + "ve changed follow the next step in your assignment task."
$find1 = $d.Content
$ok1 = $find1.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, $anchorText + "^p", 2)

# --- Step 2: append a second new paragraph (also inherits the border)
# directly after the blank one, using a real paragraph-mark character with
# Range.InsertAfter so no placeholder/empty run gets left behind in the
# blank paragraph.
$blank = $d.Paragraphs.Item(6)
$cr = [char]13
$content = "Add new content for the files testing"
$blank.Range.InsertAfter($cr + $content)

$textPara = $d.Paragraphs.Item(7)

# --- Step 3: add the "_GoBack" bookmark right after the inserted run, i.e.
# an empty/collapsed bookmark sitting between the end of the text and the
# paragraph mark. A collapsed range exactly at "end of paragraph" gets
# special-cased by this host to span the whole paragraph, so a one-off
# sentinel character is appended, the (still mid-run) collapsed bookmark
# range is created just before it, and the sentinel is deleted afterwards --
# leaving bookmarkStart/bookmarkEnd immediately after the run.
$sentinel = "~"
$paraStart = $textPara.Range.Start
$textPara.Range.InsertAfter($sentinel)

$bmPos = $paraStart + $content.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$sentinelRange = $d.Range($bmPos, $bmPos + 1)
$sentinelRange.Delete()
